# Regenerate orders with updated distance/sizes:
#   D51 -> D55, D64 -> D69, D80 -> D86, S30 -> S31
# These substitutions are applied to every text value in the used range
# (Condition, Filename_Left, Filename_Right, Distance, Size columns all
# contain these tokens as part of their strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            if ($v -like "*D51*" -or $v -like "*D64*" -or $v -like "*D80*" -or $v -like "*S30*") {
                $nv = $v.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
                if ($nv -ne $v) {
                    $cell.Value2 = $nv
                }
            }
        }
    }
}
